$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin price/volume refresh from GitHub Actions data pull.
# D-column "Price" values are textual (European-style dotted groupings,
# or decimals whose trailing zeros matter), so any cell whose new value
# could be auto-sniffed as a number is pinned to Text format first so it
# round-trips as a literal string rather than a numeric value.

$ws.Range("D2").Value = "62.070.79"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "3.415.86"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "408.88"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.08"
$ws.Range("E6").Value = "  -3.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.639"
$ws.Range("E7").Value = "  +7.81%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.733"
$ws.Range("E9").Value = "  +6.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.143"
$ws.Range("E10").Value = "  +17.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.48"
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000220"
$ws.Range("E12").Value = "  +69.19%  "
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("D14").Value = "3.960.91"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.94"
$ws.Range("E15").Value = "  +6.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.79"
$ws.Range("E16").Value = "  +4.53%  "
$ws.Range("D17").Value = "3.405.34"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.15"
$ws.Range("E18").Value = "  +10.17%  "
$ws.Range("E19").Value = "  +5.07%  "
$ws.Range("D20").Value = "62.020.49"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "414.13"
$ws.Range("E21").Value = "  +31.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "89.14"
$ws.Range("E22").Value = "  +5.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.17"
$ws.Range("E23").Value = "  -0.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.05"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.25"
$ws.Range("E25").Value = "  +2.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "33.24"
$ws.Range("E26").Value = "  +12.18%  "
$ws.Range("E27").Value = "  +7.56%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.58"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.73"
$ws.Range("E30").Value = "  -4.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.92"
$ws.Range("E31").Value = "  +4.84%  "
$ws.Range("E32").Value = "  -2.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.116"
$ws.Range("E33").Value = "  -0.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.72"
$ws.Range("E34").Value = "  +0.13%  "
$ws.Range("E35").Value = "  +0.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0500"
$ws.Range("E36").Value = "  +3.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.22"
$ws.Range("E37").Value = "  +4.82%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.36"
$ws.Range("E39").Value = "  -1.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.134"
$ws.Range("E40").Value = "  +6.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.91"
$ws.Range("E41").Value = "  -1.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.312"
$ws.Range("E42").Value = "  +4.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "141.38"
$ws.Range("E43").Value = "  +2.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.97"
$ws.Range("E44").Value = "  -0.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.10"
$ws.Range("E45").Value = "  +2.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.42"
$ws.Range("E46").Value = "  +8.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.61"
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.91"
$ws.Range("E48").Value = "  +2.65%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "3.761.11"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.113.34"
$ws.Range("E50").Value = "  -0.43%  "
$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.37"
$ws.Range("E51").Value = "  +2.53%  "
